# Update the crypto price/volume table (Sheet1) to reflect the
# Sat Dec 31 17:52:00 UTC 2022 GitHub Actions data refresh.
#
# Column D (Price) values are written with a leading apostrophe so the
# numeric-looking text ("246.56", "3.600", ...) is stored as TEXT
# (matching the workbook's inline-string cells) instead of being
# auto-converted to a number by Excel's input parser.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.56"
$ws.Range("D3").Value = "'26.59"
$ws.Range("D5").Value = "'0.05609"
$ws.Range("D6").Value = "'6.477"
$ws.Range("D7").Value = "'0.8136"
$ws.Range("D8").Value = "'0.8459"
$ws.Range("B9").Value = "BitrueCoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D9").Value = "'0.02844"
$ws.Range("E9").Value = "8BitrueCoinBTR"
$ws.Range("B10").Value = "BitMartToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D10").Value = "'0.09386"
$ws.Range("E10").Value = "9BitMartTokenBMX"
$ws.Range("B11").Value = "BitForexToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D11").Value = "'0.001513"
$ws.Range("E11").Value = "10BitForexTokenBF"
$ws.Range("B12").Value = "TigerCash"
$ws.Range("C12").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D12").Value = "'0.006132"
$ws.Range("E12").Value = "11TigerCashTCH"
$ws.Range("B13").Value = "LEO"
$ws.Range("C13").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D13").Value = "'3.600"
$ws.Range("E13").Value = "12LEOLEO"
$ws.Range("B14").Value = "GateToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D14").Value = "'3.009"
$ws.Range("E14").Value = "13GateTokenGT"
$ws.Range("B15").Value = "BTSEToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D15").Value = "'2.055"
$ws.Range("E15").Value = "14BTSETokenBTSE"
$ws.Range("B16").Value = "BitpandaEcosystemToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D16").Value = "'0.3207"
$ws.Range("E16").Value = "15BitpandaEcosystemTokenBEST"
$ws.Range("B17").Value = "WazirX"
$ws.Range("C17").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D17").Value = "'0.1341"
$ws.Range("E17").Value = "16WazirXWRX"
$ws.Range("B18").Value = "MandalaExchangeToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D18").Value = "'0.06953"
$ws.Range("E18").Value = "17MandalaExchangeTokenMDX"
$ws.Range("B19").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C19").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D19").Value = "'0.03173"
$ws.Range("E19").Value = "18LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1318"
$ws.Range("E20").Value = "19ProBitTokenPROB"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'3.745"
$ws.Range("E21").Value = "20MCDexMCB"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.04642"
$ws.Range("E22").Value = "21CoinExTokenCET"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.1350"
$ws.Range("E23").Value = "22ZBTokenZB"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.0005972"
$ws.Range("E24").Value = "23OneONE"
$ws.Range("D26").Value = "'0.004590"
$ws.Range("D40").Value = "'0.03666"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1351"
$ws.Range("E41").Value = "40BKEXTokenBKKBestin24h"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002661"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003377"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("D44").Value = "'0.008910"
$ws.Range("D45").Value = "'0.00005292"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
